$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: paragraph "Nota: El código fuente definitivo está en la
# rama………" -> the trailing ellipsis run is replaced by three runs:
#   ' "'  (space + left curly quote)
#   'main' (wrapped in proofErr spellStart/spellEnd, as Word's spell
#           checker would flag the English word inside Spanish text)
#   '".'  (right curly quote + period)
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Nota: El*rama*") {
        $rng = $p.Range
        $rng.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4AB1DEA7" w14:textId="551837F9" w:rsidR="000A529D" w:rsidRDefault="000A529D"><w:r><w:t>Nota: El c&#243;digo fuente definitivo est&#225; en la rama</w:t></w:r><w:r><w:t xml:space="preserve"> &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221;.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
        break
    }
}

# ------------------------------------------------------------------
# Change 2: paragraph made of 5 separate runs (XML Y / XSD / validando.../
# space / url) is collapsed into a single run with the full sentence.
# ------------------------------------------------------------------
$old = "Tal y como se pide en el enunciado, cada uno de los miembros del equipo ha propuesto una solución de XML Y XSD validando el documento con el validador online: http://www.xmlvalidation.com"
$found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2)
Write-Output "Change2 found/replaced: $found"
